# Proctor Provisioning code added
# Update the STAGE sheet's row 2 (the "exam taker" record) with the new
# Proctor Provisioning sample data, matching the PROD sheet's row 2 layout
# (shared formats), refresh the view/selection state on both sheets, and
# tone down the Examtaker/Location font color (theme text color instead of
# the old blue).

$wb = $excel.ActiveWorkbook

$stage = $wb.Worksheets.Item("STAGE")
$prod  = $wb.Worksheets.Item("PROD")

# --- STAGE!A2:H2 -------------------------------------------------------
$stage.Cells.Item(2, 1).Value = "december2adsff8exargmtaker01@gmail.com"
$stage.Cells.Item(2, 2).Value = "222211"
$stage.Cells.Item(2, 3).Value = "Examtaker"
$stage.Cells.Item(2, 4).Value = 4599
$stage.Cells.Item(2, 5).Value = "Chicago"
$stage.Cells.Item(2, 6).Value = 2008
$stage.Cells.Item(2, 7).Value = "JUL"
$stage.Cells.Item(2, 8).Value = 16

# Make sure F2 carries an integer number format like the PROD sheet's F2
# (C2/D2/G2/H2 are already plain "General" cells, same as the default).
$stage.Range("F2").NumberFormat = "0"

# --- PROD!E2 -------------------------------------------------------
# Text content is unchanged (still "BBSR") but keep it explicit.
$prod.Cells.Item(2, 5).Value = "BBSR"

# --- Font color tweak (Examtaker/Location 8pt Arial font) --------------
# Used by STAGE!E2 (style 7). Switch from the hard-coded blue to the
# workbook's theme text color.
$stage.Range("E2").Font.ThemeColor = 1

# --- Selections / active sheet -----------------------------------------
# PROD keeps A2:H2 selected (anchored at A2) and is no longer the active tab.
$prod.Range("A2:H2").Select()

# STAGE becomes the tab that is selected/active, with H10 selected - do
# this last so STAGE ends up as the active/visible sheet.
$stage.Activate()
$stage.Range("H10").Select()
